$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows being appended to the status table (rows 67-70).
# Column layout: A=Date(serial), B=#Total Users, C=#Logged-in Users,
# D=0 Errors, E=1 Errors, F=2 Errors, G=3-5 Errors, H=6-10 Errors, I=>10 Errors
$data = @(
    @(46031, 5594, 4161, 3838, 234, 51, 35, 3, 0),
    @(46034, 5603, 4383, 4021, 267, 51, 32, 11, 1),
    @(46035, 5603, 4436, 4057, 273, 69, 33, 4, 0),
    @(46036, 5602, 4401, 4081, 224, 55, 37, 4, 0)
)

$startRow = 67
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowValues = $data[$i]

    # Column A holds the date serial; make sure new rows (69/70) pick up the
    # same date number format ("d-mmm-yy", i.e. style s="4") already used by
    # the rest of the column. Existing rows (67/68) already carry that style.
    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Value = $rowValues[0]
    $aCell.NumberFormat = "d-mmm-yy"

    for ($c = 2; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}

# Move the sheet's active selection to the newly added last row, matching
# how the workbook previously tracked the last data row (A66:I66 -> A70:I70).
$ws.Range("A70:I70").Select()
